$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the rnaSampleNumber values in column C (rows 2-27): shift by +26
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    $cell.Value = $current + 26
}

# Update the active selection to reflect where the user ended up (C28)
$ws.Range("C28").Select()
